# Revert "P0 source datasheet updated with out_of_stock sheet, moved out of expense_reports"
#
# This undoes the addition of the out_of_stock sheet: the two header
# cells it contained ("Items Out of Stock" / "Items Not Found") move back
# onto the expense_reports sheet (F1/G1), the out_of_stock sheet itself is
# deleted, and the vendor_inventory sheet loses the trailing blank row it
# had picked up along the way.

$wb = $excel.ActiveWorkbook

# --- expense_reports: restore the two header cells that used to live here
$expense = $wb.Worksheets.Item("expense_reports")
$expense.Range("F1").Value = "Items Out of Stock"
$expense.Range("G1").Value = "Items Not Found"
$expense.Range("E15").Select()

# --- vendor_inventory: drop the stray blank row 14 and reset its selection
$vendorInv = $wb.Worksheets.Item("vendor_inventory")
$vendorInv.Rows.Item(14).Delete()

# --- remove the out_of_stock sheet entirely
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("out_of_stock").Delete()
$excel.DisplayAlerts = $true

# --- vendor_inventory ends up the active/selected tab
$vendorInv.Activate()
$vendorInv.Range("D20").Select()
